$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.6731329884640765, 0.7653800872874396, 0.9857854874064238, 0.8923015287106822, 0.3617455065250397, 0.4386407136917114, 0.02572542615234852, 0.2443276047706604)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
